$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.94080000000001

$ws.Range("B4").Value = 5.559099999999994
$ws.Range("D4").Value = -8.152299999999995
$ws.Range("E4").Value = 13.5338

$ws.Range("D5").Value = -8.582599999999999

$ws.Range("B6").Value = 9.474399999999992
$ws.Range("D6").Value = -9.220799999999988

$ws.Range("B7").Value = 5.208699999999994

$ws.Range("B8").Value = 5.201499999999993
$ws.Range("D8").Value = -7.992699999999998

$ws.Range("E9").Value = 14.33140000000001

$ws.Range("E11").Value = 13.1214

$ws.Range("E14").Value = 13.12490000000001

$ws.Range("B16").Value = 8.44380000000001
$ws.Range("D16").Value = -7.853800000000001

$ws.Range("E18").Value = 12.4829

$ws.Range("B20").Value = 5.643899999999993

$ws.Range("B21").Value = 5.301899999999993

$ws.Range("D22").Value = -7.952700000000001

$ws.Range("E25").Value = 13.14310000000001
